$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update compiler-settings label for both result rows (was "Balanced")
$ws.Range("F3").Value = "Area (Aggressive)"
$ws.Range("F4").Value = "Area (Aggressive)"

# Update the synthesis results (row 3)
$ws.Range("G3").Value = 1635
$ws.Range("H3").Value = 2235
$ws.Range("J3").Value = 129.45
$ws.Range("K3").Value = 165

# Update the synthesis results (row 4)
$ws.Range("G4").Value = 1635
$ws.Range("H4").Value = 2235
$ws.Range("J4").Value = 129.45
$ws.Range("K4").Value = 165

# Add a note about additional compiler settings below the table
$ws.Range("F6").Value = "Note: we also changed additional compiler settings (fitter effort, etc). Have included QPF files for reference."

# Update view state to match the authored workbook
$excel.ActiveWindow.Zoom = 145
$ws.Range("F15").Select()
